$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 99
$ws.Range("D99").Value = 45005
$ws.Range("J99").Value = 600
$ws.Range("K99").Value = 15000
$ws.Range("L99").Value = 16000
$ws.Range("M99").Value = 15417
$ws.Range("P99").Value = 1186

# Row 100
$ws.Range("D100").Value = 44789
$ws.Range("J100").Value = 580
$ws.Range("K100").Value = 11000
$ws.Range("L100").Value = 12000
$ws.Range("M100").Value = 11448
$ws.Range("P100").Value = 881

# Row 101
$ws.Range("D101").Value = 44873
$ws.Range("J101").Value = 280
$ws.Range("K101").Value = 12000
$ws.Range("L101").Value = 13000
$ws.Range("M101").Value = 12536
$ws.Range("P101").Value = 964

# Row 102
$ws.Range("D102").Value = 44508
$ws.Range("J102").Value = 400
$ws.Range("K102").Value = 13000
$ws.Range("L102").Value = 15000
$ws.Range("M102").Value = 13850
$ws.Range("P102").Value = 1065

# Row 103
$ws.Range("D103").Value = 44663
$ws.Range("J103").Value = 150
$ws.Range("K103").Value = 10000
$ws.Range("L103").Value = 10000
$ws.Range("M103").Value = 10000
$ws.Range("P103").Value = 769

# Row 104
$ws.Range("D104").Value = 44760
$ws.Range("K104").Value = 12000
$ws.Range("L104").Value = 13000
$ws.Range("M104").Value = 12425
$ws.Range("P104").Value = 956

# Row 105
$ws.Range("D105").Value = 44998
$ws.Range("K105").Value = 14000
$ws.Range("L105").Value = 15000
$ws.Range("M105").Value = 14425
$ws.Range("P105").Value = 1110

# Row 106
$ws.Range("D106").Value = 44740
$ws.Range("J106").Value = 400
$ws.Range("K106").Value = 13000
$ws.Range("L106").Value = 14000
$ws.Range("M106").Value = 13425
$ws.Range("P106").Value = 1033

# Row 107
$ws.Range("D107").Value = 44328
$ws.Range("J107").Value = 290
$ws.Range("K107").Value = 23000
$ws.Range("L107").Value = 25000
$ws.Range("M107").Value = 23828
$ws.Range("P107").Value = 1833

# Row 108
$ws.Range("D108").Value = 44795
$ws.Range("J108").Value = 160
$ws.Range("L108").Value = 12000
$ws.Range("M108").Value = 12000
$ws.Range("P108").Value = 923

# Row 109
$ws.Range("D109").Value = 44651
$ws.Range("J109").Value = 180
$ws.Range("K109").Value = 12000
$ws.Range("L109").Value = 14000
$ws.Range("M109").Value = 12889
$ws.Range("P109").Value = 991

# Row 110
$ws.Range("I110").Value = "Primera"
$ws.Range("J110").Value = 140
$ws.Range("K110").Value = 14000
$ws.Range("L110").Value = 15000
$ws.Range("M110").Value = 14571
$ws.Range("P110").Value = 1121

# Row 111
$ws.Range("D111").Value = 44571
$ws.Range("I111").Value = "Segunda"
$ws.Range("J111").Value = 30
$ws.Range("K111").Value = 12000
$ws.Range("M111").Value = 12000
$ws.Range("P111").Value = 923

# Row 112
$ws.Range("D112").Value = 44860
$ws.Range("J112").Value = 580
$ws.Range("K112").Value = 11000
$ws.Range("L112").Value = 12000
$ws.Range("M112").Value = 11448
$ws.Range("P112").Value = 881

# Row 113
$ws.Range("D113").Value = 44473
$ws.Range("J113").Value = 250
$ws.Range("K113").Value = 12000
$ws.Range("L113").Value = 13000
$ws.Range("M113").Value = 12400
$ws.Range("P113").Value = 954

# Row 114
$ws.Range("D114").Value = 44306
$ws.Range("J114").Value = 230
$ws.Range("K114").Value = 24000
$ws.Range("L114").Value = 25000
$ws.Range("M114").Value = 24435
$ws.Range("P114").Value = 1880

# Row 115
$ws.Range("D115").Value = 44168
$ws.Range("J115").Value = 500
$ws.Range("K115").Value = 30000
$ws.Range("L115").Value = 32000
$ws.Range("M115").Value = 31080
$ws.Range("P115").Value = 2391
